# StructureDefinition-ror-organization-additional-name.xlsx
# 1. Metadata!B8 - bump the "Date" value to the new publication timestamp.
# 2. Elements sheet - the two right-most "Mapping" columns (AK = RIM Mapping,
#    AL = Spécification métier) were reordered so the business-mapping
#    column now comes first; swap their header/data content and widths.

$wb = $excel.ActiveWorkbook

# --- 1. Metadata date -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# --- 2. Elements sheet: swap columns AK (37) and AL (38) --------------
$ws = $wb.Worksheets.Item("Elements")

for ($r = 1; $r -le 6; $r++) {
    $akAddr = "AK" + $r
    $alAddr = "AL" + $r
    $akVal = $ws.Range($akAddr).Value2
    $alVal = $ws.Range($alAddr).Value2
    $ws.Range($akAddr).Value = $alVal
    $ws.Range($alAddr).Value = $akVal
}

# Column widths also swap places: AK was the narrow "RIM Mapping" column
# (~24.98 chars) and AL was the wide "Spécification métier" column
# (~71.57 chars); after the move AK is wide and AL is narrow.
$ws.Columns.Item(37).ColumnWidth = 71.5703125 - 0.8333333333333334
$ws.Columns.Item(38).ColumnWidth = 24.98046875 - 0.8333333333333334
